$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the handful of D-column cells whose new value would otherwise be
# auto-coerced to a number (losing a significant trailing zero) to stay text,
# matching the "General -> stays text" behaviour of the other price cells.
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "91.630.61"
$ws.Range("E2").Value = "  +1.17%  "

# Row 3
$ws.Range("D3").Value = "3.155.06"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "240.32"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("D6").Value = "619.83"
$ws.Range("E6").Value = "  -0.98%  "

# Row 7
$ws.Range("D7").Value = "1.12"
$ws.Range("E7").Value = "  -0.85%  "

# Row 8
$ws.Range("D8").Value = "0.388"
$ws.Range("E8").Value = "  +4.20%  "

# Row 9
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").Value = "3.154.45"
$ws.Range("E10").Value = "  +14.64%  "

# Row 11
$ws.Range("D11").Value = "0.744"
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
$ws.Range("E12").Value = "  +0.44%  "

# Row 13
$ws.Range("E13").Value = "  +2.01%  "

# Row 14
$ws.Range("D14").Value = "35.11"
$ws.Range("E14").Value = "  -0.51%  "

# Row 15
$ws.Range("D15").Value = "5.60"
$ws.Range("E15").Value = "  +1.98%  "

# Row 16
$ws.Range("D16").Value = "91.350.02"
$ws.Range("E16").Value = "  +0.97%  "

# Row 18
$ws.Range("D18").Value = "3.152.08"
$ws.Range("E18").Value = "  +0.06%  "

# Row 19
$ws.Range("D19").Value = "3.73"
$ws.Range("E19").Value = "  -2.83%  "

# Row 20
$ws.Range("D20").Value = "15.00"
$ws.Range("E20").Value = "  +5.31%  "

# Row 21
$ws.Range("D21").Value = "5.91"
$ws.Range("E21").Value = "  +2.25%  "

# Row 22
$ws.Range("D22").Value = "457.07"
$ws.Range("E22").Value = "  +2.54%  "

# Row 23
$ws.Range("D23").Value = "0.0000203"
$ws.Range("E23").Value = "  -3.19%  "

# Row 24
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  +1.20%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "5.91"
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "1.62"
$ws.Range("E26").Value = "  +61.57%  "

# Row 27
$ws.Range("D27").Value = "88.66"
$ws.Range("E27").Value = "  -4.66%  "

# Row 28
$ws.Range("D28").Value = "11.81"
$ws.Range("E28").Value = "  -2.04%  "

# Row 29
$ws.Range("D29").Value = "3.317.67"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("D30").Value = "0.148"
$ws.Range("E30").Value = "  +36.13%  "

# Row 31
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("D32").Value = "0.230"
$ws.Range("E32").Value = "  +5.49%  "

# Row 33
$ws.Range("E33").Value = "  -4.89%  "

# Row 34
$ws.Range("D34").Value = "9.38"
$ws.Range("E34").Value = "  +1.17%  "

# Row 35
$ws.Range("E35").Value = "  +11.98%  "

# Row 36
$ws.Range("D36").Value = "26.38"
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("D37").Value = "7.48"
$ws.Range("E37").Value = "  -1.12%  "

# Row 38
$ws.Range("D38").Value = "1.95"
$ws.Range("E38").Value = "  +1.42%  "

# Row 39
$ws.Range("D39").Value = "493.26"
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("B40").Value = "MantraDAO"
$ws.Range("C40").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D40").Value = "3.89"
$ws.Range("E40").Value = "  -12.07%  "

# Row 41
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "1.32"
$ws.Range("E41").Value = "  +2.34%  "

# Row 42
$ws.Range("D42").Value = "0.443"
$ws.Range("E42").Value = "  +6.26%  "

# Row 43
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -6.17%  "

# Row 44
$ws.Range("D44").Value = "22.16"
$ws.Range("E44").Value = "  +0.17%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.94"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.709"
$ws.Range("E47").Value = "  +3.37%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "155.72"
$ws.Range("E48").Value = "  -2.26%  "

# Row 49
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").Value = "  +1.47%  "

# Row 50
$ws.Range("D50").Value = "4.48"
$ws.Range("E50").Value = "  -1.74%  "

# Row 51
$ws.Range("D51").Value = "44.11"
